# "convocazioni chiuse colosseo agosto"
# Mark every member on the "data" sheet as no longer active: column D
# ("Attivo") flips from "si" to "no" for all data rows (2-142).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# Bulk-update the whole "Attivo" column in one shot.
$ws.Range("D2:D142").Value = "no"

# Restore the view/selection state left behind by the edit.
$ws.Range("C130").Select()
